# Test interlink for hybrid ac-dc systems
$wb = $excel.ActiveWorkbook

# --- Bus sheet: selection moves from K10 to C5 ---
$wsBus = $wb.Worksheets.Item("Bus")
$wsBus.Select()
$wsBus.Range("C5").Select()

# --- Device sheet: clear the test values in C5:F5, selection moves to B5 ---
$wsDevice = $wb.Worksheets.Item("Device")
$wsDevice.Select()
$wsDevice.Range("C5:F5").ClearContents()
$wsDevice.Range("B5").Select()

# --- Advance sheet: no longer the active tab, selection stays at B12 ---
$wsAdvance = $wb.Worksheets.Item("Advance")
$wsAdvance.Select()
$wsAdvance.Range("B12").Select()

# --- Basic sheet: becomes the active tab, selection moves to B10 ---
$wsBasic = $wb.Worksheets.Item("Basic")
$wsBasic.Select()
$wsBasic.Range("B10").Select()
